$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lab_abstract")
Write-Host $ws.Name
